$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "41.642.79"
$ws.Range("E2").Value = "  -1.52%  "

# Row 3
$ws.Range("D3").Value = "2.171.69"
$ws.Range("E3").Value = "  -2.75%  "

# Row 4
$ws.Range("E4").Value = "  -0.13%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "238.47"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.02%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "72.24"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.17%  "

# Row 8
$ws.Range("E8").Value = "  -0.14%  "

# Row 9
$ws.Range("E9").Value = "  -4.82%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.96"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.65%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0909"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.52%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.45"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.00%  "

# Row 13
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.72"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.07%  "

# Row 14
$ws.Range("B14").Value = "TRON"
$ws.Range("C14").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0998"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.81%  "

# Row 15
$ws.Range("D15").Value = "2.496.70"
$ws.Range("E15").Value = "  -2.87%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.31"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.26%  "

# Row 17
$ws.Range("D17").Value = "2.156.24"
$ws.Range("E17").Value = "  -4.23%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.781"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -7.06%  "

# Row 19
$ws.Range("D19").Value = "41.506.65"
$ws.Range("E19").Value = "  -1.57%  "

# Row 20
$ws.Range("E20").Value = "  -2.90%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "70.11"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.09%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.79"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -7.19%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.77"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -13.90%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "226.78"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.84%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.03"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.50%  "

# Row 26
$ws.Range("E26").Value = "  +0.30%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.75"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.23%  "

# Row 28
$ws.Range("E28").Value = "  -9.74%  "

# Row 29
$ws.Range("E29").Value = "  -3.73%  "

# Row 30
$ws.Range("E30").Value = "  -1.59%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "171.03"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.49%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.84"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.80%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "33.61"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +11.19%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0774"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.58%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.24"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -7.64%  "

# Row 36
$ws.Range("E36").Value = "  -3.63%  "

# Row 37
$ws.Range("E37").Value = "  -1.53%  "

# Row 38
$ws.Range("E38").Value = "  -4.94%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0307"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.36%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "12.14"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -8.65%  "

# Row 41
$ws.Range("E41").Value = "  -1.98%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.38"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.71%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "58.97"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -9.47%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.45"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.26%  "

# Row 45
$ws.Range("E45").Value = "  -5.06%  "

# Row 46
$ws.Range("E46").Value = "  -3.85%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "97.53"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.85%  "

# Row 48
$ws.Range("E48").Value = "  -3.64%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.11"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.72%  "

# Row 50
$ws.Range("E50").Value = "  -7.57%  "

# Row 51
$ws.Range("E51").Value = "  -2.38%  "
